$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Cu3Au", 0.633, 40, -3, 38),
    @("Cu3Au", 0.58, 54, -6, 48),
    @("Cu3Au", 0.535, 57, -2, 55),
    @("Cu3Au", 0.486, 77, -15, 62),
    @("Cu3Au", 0.43, 86, 11, 75),
    @("NaZn13", 0.67, 18, -5, 14),
    @("NaZn13", 0.658, 14, -10, 5),
    @("NaZn13", 0.633, 10, -8, 3),
    @("NaZn13", 0.58, 10, -10, 0)
)

$startRow = 26
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}

$ws.Range("H31").Select()
